$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -6.907599999999997
$ws.Range("D6").Value = -8.155100000000003
$ws.Range("D7").Value = -7.519099999999994
$ws.Range("D16").Value = -8.1297
$ws.Range("D20").Value = -8.318100000000001
